$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 3099.503889238888
$ws.Range("D2").Value = 117845141.8069585
$ws.Range("E2").Value = 9353990175.932438
$ws.Range("G2").Value = 9471838417.918615
